$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Beetroot" row (row 6), which only had a partial record
# (no Voc/Isc/FF data). This shifts rows 7-19 up by one.
$ws.Rows.Item(6).Delete()

# Rename the "*-simp" dye labels to their final names (rows 15-18
# after the shift above). Assign starting from the bottom row so the
# new shared-string table entries are appended in T4,T3,T2,T1 order.
$ws.Range("A18").Value = "T4"
$ws.Range("A17").Value = "T3"
$ws.Range("A16").Value = "T2"
$ws.Range("A15").Value = "T1"

# Update the sheet view: Excel had scrolled so that row 4 was the top
# visible row and H5 was selected; reset that back to the top of the
# sheet with A16 selected (as saved with the new, shorter dataset).
$ws.Range("A16").Select()
$excel.ActiveWindow.ScrollRow = 1
